# Insert a new data row at Excel row 411 (pushing the existing rows 411-434
# down to 412-435) and populate it with a new "Pepino ensalada" price-report
# entry for Vega Modelo de Temuco, Región del Maule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 411..434 down to 412..435, leaving a blank row 411 behind
# (formatting is inherited from the row above, matching column D's date style).
$ws.Rows(411).Insert()

$ws.Range("A411").Value = 10
$ws.Range("B411").Value = "Vega Modelo de Temuco"
$ws.Range("C411").Value = "La Araucanía"
$ws.Range("D411").Value = 44610
$ws.Range("E411").Value = 9
$ws.Range("F411").Value = 100112043
$ws.Range("G411").Value = "Pepino ensalada"
$ws.Range("H411").Value = "Sin especificar"
$ws.Range("I411").Value = "Extra"
$ws.Range("J411").Value = 30
$ws.Range("K411").Value = 17000
$ws.Range("L411").Value = 17000
$ws.Range("M411").Value = 17000
$ws.Range("N411").Value = '$/caja 50 unidades'
$ws.Range("O411").Value = "Región del Maule"
$ws.Range("P411").Value = 340
$ws.Range("Q411").Value = 50
$ws.Range("R411").Value = "Hortaliza"
